# MCP3561 SCAN Mode Samplerate Calculator
# Commit: "SCAN mode working, 32bit format unclear"
#
# The author changed the two primary inputs of the calculator:
#   N_CH (C9): number of scanned channels   1 -> 4
#   OSR  (C10): oversampling ratio          98304 -> 4096
# Every other changed cell on the sheet is a formula that depends
# (directly or transitively) on C9/C10, so it is expected to recompute
# automatically once the inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# N_CH: number of channels scanned, 1 -> 4
$ws.Range("C9").Value = 4

# OSR: oversampling ratio, 98304 -> 4096
$ws.Range("C10").Value = 4096

# Force a full recalculation so all dependent formula cells (C7, F7, F10,
# F11, C12, F12, C16, F16, C17, C18, C20, C21, ...) refresh their cached
# <v> results.
$excel.CalculateFullRebuild()
$excel.Calculate()

# The author's cursor ended up on C21 (the final "effective sampling
# rate" result) when the workbook was saved.
$ws.Range("C21").Select()

# Best-effort: also mirror the saved window geometry from the diff.
try {
    $win = $excel.ActiveWindow
    $win.WindowState = -4143  # xlNormal
    $win.Left = 5670
    $win.Top = 4680
    $win.Width = 37965
    $win.Height = 15885
} catch {
    # Window placement isn't part of the workbook's data model in every
    # host; ignore if unsupported.
}
